$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "B08L7BD4C6"

# Update column A values (A1:A100) with the new product-title strings
$values = @(
    "black sport",
    "clothing set",
    "seamless legging",
    "seamless crop top",
    "legging suit",
    "legging top set",
    "green apparel",
    "small crop top",
    "black compression",
    "winter activewear",
    "winter blue",
    "winter tracksuit",
    "pink work shirt",
    "work out outfit",
    "winter outfit",
    "2 piece matching set",
    "2 piece long sleeve",
    "2 piece workout set",
    "black pant suit",
    "black pant set",
    "black workout shirt",
    "black workout top",
    "seamless top",
    "pink tight shirt",
    "pink crop shirt",
    "green high top",
    "work out crop",
    "two piece legging set",
    "seamless workout top",
    "seamless workout",
    "black sport shirt",
    "winter workout",
    "long black top",
    "tight long sleeve workout shirt",
    "workout tracksuit",
    "workout pant",
    "2 piece jogger",
    "2 piece woman",
    "athletic legging",
    "black tight top",
    "black jogger",
    "black long crop top",
    "crop top woman",
    "green black top",
    "long sleeve winter shirt",
    "2 piece pant set",
    "2 piece pant suit",
    "long sleeve work out shirt",
    "long sleeve set",
    "black top woman",
    "long sleeve outfit",
    "crop top pant set",
    "long sleeve athletic top",
    "long sleeve two piece set",
    "pink workout set",
    "black jogger set",
    "tight workout top",
    "long sleeve clothing",
    "black athletic crop top",
    "small black crop top",
    "black athletic sleeve",
    "matching legging set",
    "long sleeve workout top",
    "long sleeve crop top workout",
    "pink long sleeve crop",
    "black long sleeve crop",
    "green long sleeve crop",
    "long sleeve crop top set",
    "long sleeve work",
    "top pant set",
    "pink outfit set",
    "black crop",
    "long black",
    "work outfit",
    "piece work",
    "set two",
    "fitness pant",
    "winter shirt",
    "workout tight",
    "high top woman",
    "matching pant set",
    "gym tight",
    "legging outfit",
    "long green",
    "crop legging",
    "work set",
    "seamless long sleeve workout top",
    "two apparel",
    "tracksuit pant",
    "work out pink",
    "workout track",
    "black high top",
    "fashionable top",
    "winter top",
    "2 small",
    "seamless crop",
    "yoga with",
    "pink out shirt",
    "pink crop",
    "yoga jogger"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

# Update the active selection to E13
$ws.Range("E13").Select()

